$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(4, 6).Value = 4836
$ws1.Cells.Item(5, 6).Value = 213
$ws1.Cells.Item(6, 6).Value = 161
$ws1.Cells.Item(7, 6).Value = 123
$ws1.Cells.Item(11, 6).Value = 0
$ws1.Cells.Item(12, 6).Value = 1200
$ws1.Cells.Item(13, 6).Value = 0
$ws1.Cells.Item(14, 6).Value = 0
$ws1.Cells.Item(15, 6).Value = 190
$ws1.Cells.Item(17, 6).Value = 1
$ws1.Cells.Item(18, 6).Value = 154
$ws1.Cells.Item(19, 6).Value = 114
$ws1.Cells.Item(20, 6).Value = 0
$ws1.Cells.Item(23, 6).Value = 39
$ws1.Cells.Item(24, 6).Value = 88
$ws1.Cells.Item(25, 6).Value = 0
$ws1.Cells.Item(27, 6).Value = 3991
$ws1.Cells.Item(31, 6).Value = 2592
$ws1.Cells.Item(33, 6).Value = 534
$ws1.Cells.Item(36, 6).Value = 314
$ws1.Cells.Item(38, 6).Value = 182
$ws1.Cells.Item(40, 6).Value = 1571
$ws1.Cells.Item(43, 6).Value = 76
$ws1.Cells.Item(44, 6).Value = 0
$ws1.Cells.Item(45, 6).Value = 502
$ws1.Cells.Item(46, 6).Value = 483
$ws1.Cells.Item(47, 6).Value = 2

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 112
$ws4.Cells.Item(3, 6).Value = 234
$ws4.Cells.Item(4, 6).Value = 4836
$ws4.Cells.Item(5, 6).Value = 213
$ws4.Cells.Item(7, 6).Value = 123
$ws4.Cells.Item(8, 6).Value = 0
$ws4.Cells.Item(9, 6).Value = 0
$ws4.Cells.Item(11, 6).Value = 765
$ws4.Cells.Item(12, 6).Value = 229
$ws4.Cells.Item(13, 6).Value = 1200
$ws4.Cells.Item(14, 6).Value = 0
$ws4.Cells.Item(15, 6).Value = 190
$ws4.Cells.Item(17, 6).Value = 1
$ws4.Cells.Item(18, 6).Value = 154
$ws4.Cells.Item(19, 6).Value = 114
$ws4.Cells.Item(20, 6).Value = 4066
$ws4.Cells.Item(21, 6).Value = 6378
$ws4.Cells.Item(23, 6).Value = 0
$ws4.Cells.Item(25, 6).Value = 0
$ws4.Cells.Item(26, 6).Value = 0
$ws4.Cells.Item(27, 6).Value = 0
$ws4.Cells.Item(28, 6).Value = 409
$ws4.Cells.Item(30, 6).Value = 0
$ws4.Cells.Item(32, 6).Value = 569
$ws4.Cells.Item(35, 6).Value = 301
$ws4.Cells.Item(36, 6).Value = 314
$ws4.Cells.Item(38, 6).Value = 182
$ws4.Cells.Item(39, 6).Value = 0
$ws4.Cells.Item(40, 6).Value = 1571
$ws4.Cells.Item(41, 6).Value = 976
$ws4.Cells.Item(43, 6).Value = 0
$ws4.Cells.Item(44, 6).Value = 0
$ws4.Cells.Item(46, 6).Value = 483
$ws4.Cells.Item(47, 6).Value = 2
